# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the a2b6ad19-... item in the zh-cn and de-de sheets, and
# records the error detail message + widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 on every sheet corresponds to the a2b6ad19-0849-4c34-b20f-f7a40de40291
# handoff item. Its status flips from "Ready for handoff" to a failure state.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch in the Error Detail column
# (column P) for the a2b6ad19 row on each locale sheet.
$zhcn.Range("P3").Value = "Handback file name: sntpivxw.lbt is different with handoff file name: a2b6ad19-0849-4c34-b20f-f7a40de40291.7ae44d30cc9053fb74e76c5ceb674a369ced60a5.zh-cn."
$dede.Range("P3").Value = "Handback file name: sntpivxw.lbt is different with handoff file name: a2b6ad19-0849-4c34-b20f-f7a40de40291.7ae44d30cc9053fb74e76c5ceb674a369ced60a5.de-de."

# Widen the Error Detail column (P) so the long message is readable.
# (COM ColumnWidth is character-width units with a 5px/MDW padding baked in;
# 39.17 round-trips to an OOXML stored "width" of exactly 40, matching the
# other full-width columns already in this sheet.)
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
